# Insert a new record as the new row 3 (pushing the existing rows 3-34
# down to rows 4-35), matching the data that was added when re-uploading
# the contact list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3..34 down to 4..35 by inserting a fresh row at 3.
$ws.Rows.Item(3).Insert()

# Populate the new row. Values are entered with a leading apostrophe so
# Excel stores them as literal text (keeping the "+" prefix on phone
# numbers and the ISO date string instead of a date serial number).
$ws.Range("A3").Value = "'+5521985096467"
$ws.Range("B3").Value = "'21"
$ws.Range("C3").Value = "'2024-10-31"

# Copy the formatting/style from the row below (the original row 3, now
# row 4) onto the new row so it matches the rest of the table exactly.
$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
